$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# Row 98
$ws.Range("H98").Value = 1613.7609
$ws.Range("I98").Value = 1650.973
$ws.Range("K98").Value = 1650.973
$ws.Range("M98").Value = -152.973

# Row 112
$ws.Range("H112").Value = 912362.2
$ws.Range("I112").Value = 2099
$ws.Range("J112").Value = 1114642.9
$ws.Range("K112").Value = 6297
$ws.Range("L112").Value = 3343928.7
$ws.Range("M112").Value = -5189
$ws.Range("N112").Value = -3346144.7

# Row 122
$ws.Range("H122").Value = 1613.7609
$ws.Range("I122").Value = 1650.973
$ws.Range("K122").Value = 4952.919
$ws.Range("M122").Value = -2502.919

# Row 137
$ws.Range("H137").Value = 2891.2727
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2891.2727
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 8673.8181
$ws.Range("N137").Value = -13773.8181
$ws.Range("M137").ClearContents()


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2951.3
$ws.Range("I32").Value = 2932.7017
$ws.Range("K32").Value = 2932.7017
$ws.Range("M32").Value = -2645.7017

# Row 61
$ws.Range("H61").Value = 6858.364
$ws.Range("I61").Value = 8914.23
$ws.Range("J61").Value = 3888.7778
$ws.Range("K61").Value = 8914.23
$ws.Range("L61").Value = 3888.7778
$ws.Range("M61").Value = -8702.23
$ws.Range("N61").Value = -4312.7778

# Row 74
$ws.Range("H74").Value = 7874.5
$ws.Range("I74").Value = 19998
$ws.Range("J74").Value = 3833.3333
$ws.Range("K74").Value = 19998
$ws.Range("L74").Value = 3833.3333
$ws.Range("M74").Value = -19124
$ws.Range("N74").Value = -5581.3333

# Row 77
$ws.Range("H77").Value = 7874.5
$ws.Range("I77").Value = 19998
$ws.Range("J77").Value = 3833.3333
$ws.Range("K77").Value = 99990
$ws.Range("L77").Value = 19166.6665
$ws.Range("M77").Value = -95622
$ws.Range("N77").Value = -27902.6665

# Row 97
$ws.Range("H97").Value = 4973.933
$ws.Range("I97").Value = 2398.9583
$ws.Range("K97").Value = 2398.9583
$ws.Range("M97").Value = -1902.9583

# Row 122
$ws.Range("H122").Value = 2310.9707
$ws.Range("I122").Value = 2322.276
$ws.Range("J122").Value = 2245.4
$ws.Range("K122").Value = 6966.828
$ws.Range("L122").Value = 6736.200000000001
$ws.Range("M122").Value = -4516.828
$ws.Range("N122").Value = -11636.2

# Row 136
$ws.Range("H136").Value = 6858.364
$ws.Range("I136").Value = 8914.23
$ws.Range("J136").Value = 3888.7778
$ws.Range("K136").Value = 26742.69
$ws.Range("L136").Value = 11666.3334
$ws.Range("M136").Value = -24192.69
$ws.Range("N136").Value = -16766.3334


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5558.804
$ws.Range("I31").Value = 4080.7856
$ws.Range("J31").Value = 6205.4375
$ws.Range("K31").Value = 4080.7856
$ws.Range("L31").Value = 6205.4375
$ws.Range("M31").Value = -3785.7856
$ws.Range("N31").Value = -6795.4375

# Row 34
$ws.Range("H34").Value = 5558.804
$ws.Range("I34").Value = 4080.7856
$ws.Range("J34").Value = 6205.4375
$ws.Range("K34").Value = 4080.7856
$ws.Range("L34").Value = 6205.4375
$ws.Range("M34").Value = -3878.7856
$ws.Range("N34").Value = -6609.4375

# Row 80
$ws.Range("H80").Value = 29999
$ws.Range("J80").Value = 29999
$ws.Range("L80").Value = 29999
$ws.Range("N80").Value = -32245

# Row 83
$ws.Range("H83").Value = 29999
$ws.Range("J83").Value = 29999
$ws.Range("L83").Value = 89997
$ws.Range("N83").Value = -101229


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 6306.3
$ws.Range("J68").Value = 7445.375
$ws.Range("L68").Value = 22336.125
$ws.Range("N68").Value = -23958.125

# Row 71
$ws.Range("H71").Value = 6306.3
$ws.Range("J71").Value = 7445.375
$ws.Range("L71").Value = 67008.375
$ws.Range("N71").Value = -75120.375


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5793.1816
$ws.Range("I80").Value = 3372.7
$ws.Range("J80").Value = 29998
$ws.Range("K80").Value = 3372.7
$ws.Range("L80").Value = 29998
$ws.Range("M80").Value = -2374.7
$ws.Range("N80").Value = -31994

# Row 83
$ws.Range("H83").Value = 5793.1816
$ws.Range("I83").Value = 3372.7
$ws.Range("J83").Value = 29998
$ws.Range("K83").Value = 16863.5
$ws.Range("L83").Value = 149990
$ws.Range("M83").Value = -11871.5
$ws.Range("N83").Value = -159974

# Row 95
$ws.Range("H95").Value = 47498.25
$ws.Range("J95").Value = 47498.25
$ws.Range("L95").Value = 47498.25
$ws.Range("N95").Value = -52990.25


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3303.1667
$ws.Range("I22").Value = 2091.1428
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 2091.1428
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -1796.1428
$ws.Range("N22").Value = -5590

# Row 27
$ws.Range("H27").Value = 3303.1667
$ws.Range("I27").Value = 2091.1428
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 2091.1428
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1984.1428
$ws.Range("N27").Value = -5214

# Row 41
$ws.Range("H41").Value = 966.6667
$ws.Range("I41").Value = 966.6667
$ws.Range("K41").Value = 966.6667
$ws.Range("M41").Value = -528.6667

# Row 82
$ws.Range("H82").Value = 12804.333
$ws.Range("I82").Value = 17955.5
$ws.Range("J82").Value = 2502
$ws.Range("K82").Value = 17955.5
$ws.Range("L82").Value = 2502
$ws.Range("M82").Value = -17594.5
$ws.Range("N82").Value = -3224

# Row 85
$ws.Range("H85").Value = 12804.333
$ws.Range("I85").Value = 17955.5
$ws.Range("J85").Value = 2502
$ws.Range("K85").Value = 17955.5
$ws.Range("L85").Value = 2502
$ws.Range("M85").Value = -16707.5
$ws.Range("N85").Value = -4998

# Row 132
$ws.Range("H132").Value = 4886.8823
$ws.Range("I132").Value = 4657.6665
$ws.Range("J132").Value = 5437
$ws.Range("K132").Value = 13972.9995
$ws.Range("L132").Value = 16311
$ws.Range("M132").Value = -11442.9995
$ws.Range("N132").Value = -21371

# Row 136
$ws.Range("H136").Value = 5049.7812
$ws.Range("I136").Value = 5063.72
$ws.Range("K136").Value = 15191.16
$ws.Range("M136").Value = -12641.16

